$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$s = $c.Style
$c.Value = "'305.54"
$c.Style = $s

$c = $ws.Range("D3")
$s = $c.Style
$c.Value = "'36.04"
$c.Style = $s

$c = $ws.Range("E3")
$s = $c.Style
$c.Value = "'-3.83%"
$c.Style = $s

$c = $ws.Range("D4")
$s = $c.Style
$c.Value = "'5.119"
$c.Style = $s

$c = $ws.Range("E4")
$s = $c.Style
$c.Value = "'2.25%"
$c.Style = $s

$c = $ws.Range("D5")
$s = $c.Style
$c.Value = "'0.07872"
$c.Style = $s

$c = $ws.Range("E5")
$s = $c.Style
$c.Value = "'0.12%"
$c.Style = $s

$c = $ws.Range("D6")
$s = $c.Style
$c.Value = "'2.182"
$c.Style = $s

$c = $ws.Range("E6")
$s = $c.Style
$c.Value = "'-3.02%"
$c.Style = $s

$c = $ws.Range("D7")
$s = $c.Style
$c.Value = "'7.924"
$c.Style = $s

$c = $ws.Range("E7")
$s = $c.Style
$c.Value = "'-1.23%"
$c.Style = $s

$c = $ws.Range("B8")
$s = $c.Style
$c.Value = "'GateToken"
$c.Style = $s

$c = $ws.Range("C8")
$s = $c.Style
$c.Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$c.Style = $s

$c = $ws.Range("D8")
$s = $c.Style
$c.Value = "'4.100"
$c.Style = $s

$c = $ws.Range("E8")
$s = $c.Style
$c.Value = "'2.19%"
$c.Style = $s

$c = $ws.Range("B9")
$s = $c.Style
$c.Value = "'MXToken"
$c.Style = $s

$c = $ws.Range("C9")
$s = $c.Style
$c.Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c.Style = $s

$c = $ws.Range("D9")
$s = $c.Style
$c.Value = "'0.9203"
$c.Style = $s

$c = $ws.Range("E9")
$s = $c.Style
$c.Value = "'1.15%"
$c.Style = $s

$c = $ws.Range("B10")
$s = $c.Style
$c.Value = "'LiechtensteinCryptoassetsExchange"
$c.Style = $s

$c = $ws.Range("C10")
$s = $c.Style
$c.Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$c.Style = $s

$c = $ws.Range("D10")
$s = $c.Style
$c.Value = "'0.09690"
$c.Style = $s

$c = $ws.Range("E10")
$s = $c.Style
$c.Value = "'5.21%"
$c.Style = $s

$c = $ws.Range("B11")
$s = $c.Style
$c.Value = "'WazirX"
$c.Style = $s

$c = $ws.Range("C11")
$s = $c.Style
$c.Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$c.Style = $s

$c = $ws.Range("D11")
$s = $c.Style
$c.Value = "'0.1871"
$c.Style = $s

$c = $ws.Range("E11")
$s = $c.Style
$c.Value = "'-0.27%"
$c.Style = $s

$c = $ws.Range("B12")
$s = $c.Style
$c.Value = "'MandalaExchangeToken"
$c.Style = $s

$c = $ws.Range("C12")
$s = $c.Style
$c.Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$c.Style = $s

$c = $ws.Range("D12")
$s = $c.Style
$c.Value = "'0.08643"
$c.Style = $s

$c = $ws.Range("E12")
$s = $c.Style
$c.Value = "'1.41%"
$c.Style = $s

$c = $ws.Range("B13")
$s = $c.Style
$c.Value = "'BitrueCoin"
$c.Style = $s

$c = $ws.Range("C13")
$s = $c.Style
$c.Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$c.Style = $s

$c = $ws.Range("D13")
$s = $c.Style
$c.Value = "'0.03488"
$c.Style = $s

$c = $ws.Range("E13")
$s = $c.Style
$c.Value = "'-1.27%"
$c.Style = $s

$c = $ws.Range("B14")
$s = $c.Style
$c.Value = "'BitMartToken"
$c.Style = $s

$c = $ws.Range("C14")
$s = $c.Style
$c.Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$c.Style = $s

$c = $ws.Range("D14")
$s = $c.Style
$c.Value = "'0.09933"
$c.Style = $s

$c = $ws.Range("E14")
$s = $c.Style
$c.Value = "'-0.04%"
$c.Style = $s

$c = $ws.Range("B15")
$s = $c.Style
$c.Value = "'BitForexToken"
$c.Style = $s

$c = $ws.Range("C15")
$s = $c.Style
$c.Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$c.Style = $s

$c = $ws.Range("D15")
$s = $c.Style
$c.Value = "'0.001433"
$c.Style = $s

$c = $ws.Range("E15")
$s = $c.Style
$c.Value = "'-3.87%"
$c.Style = $s

$c = $ws.Range("B16")
$s = $c.Style
$c.Value = "'TigerCash"
$c.Style = $s

$c = $ws.Range("C16")
$s = $c.Style
$c.Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$c.Style = $s

$c = $ws.Range("D16")
$s = $c.Style
$c.Value = "'0.005676"
$c.Style = $s

$c = $ws.Range("E16")
$s = $c.Style
$c.Value = "'0.94%"
$c.Style = $s

$c = $ws.Range("B17")
$s = $c.Style
$c.Value = "'LEO"
$c.Style = $s

$c = $ws.Range("C17")
$s = $c.Style
$c.Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c.Style = $s

$c = $ws.Range("D17")
$s = $c.Style
$c.Value = "'3.461"
$c.Style = $s

$c = $ws.Range("E17")
$s = $c.Style
$c.Value = "'-0.28%"
$c.Style = $s

$c = $ws.Range("D18")
$s = $c.Style
$c.Value = "'2.467"
$c.Style = $s

$c = $ws.Range("E18")
$s = $c.Style
$c.Value = "'14.37%"
$c.Style = $s

$c = $ws.Range("D19")
$s = $c.Style
$c.Value = "'0.3425"
$c.Style = $s

$c = $ws.Range("E19")
$s = $c.Style
$c.Value = "'-1.12%"
$c.Style = $s

$c = $ws.Range("E20")
$s = $c.Style
$c.Value = "'-0.53%"
$c.Style = $s

$c = $ws.Range("D21")
$s = $c.Style
$c.Value = "'4.848"
$c.Style = $s

$c = $ws.Range("E21")
$s = $c.Style
$c.Value = "'2.17%"
$c.Style = $s

$c = $ws.Range("E22")
$s = $c.Style
$c.Value = "'0.03%"
$c.Style = $s

$c = $ws.Range("E23")
$s = $c.Style
$c.Value = "'-2.03%"
$c.Style = $s

$c = $ws.Range("D24")
$s = $c.Style
$c.Value = "'0.005091"
$c.Style = $s

$c = $ws.Range("E24")
$s = $c.Style
$c.Value = "'14.39%"
$c.Style = $s

$c = $ws.Range("D25")
$s = $c.Style
$c.Value = "'0.001232"
$c.Style = $s

$c = $ws.Range("E25")
$s = $c.Style
$c.Value = "'0.42%"
$c.Style = $s

$c = $ws.Range("D26")
$s = $c.Style
$c.Value = "'0.0001401"
$c.Style = $s

$c = $ws.Range("E26")
$s = $c.Style
$c.Value = "'7.94%"
$c.Style = $s

$c = $ws.Range("D39")
$s = $c.Style
$c.Value = "'0.01848"
$c.Style = $s

$c = $ws.Range("E39")
$s = $c.Style
$c.Value = "'5.01%"
$c.Style = $s

$c = $ws.Range("D40")
$s = $c.Style
$c.Value = "'0.04786"
$c.Style = $s

$c = $ws.Range("E40")
$s = $c.Style
$c.Value = "'1.14%"
$c.Style = $s

$c = $ws.Range("D41")
$s = $c.Style
$c.Value = "'0.007746"
$c.Style = $s

$c = $ws.Range("E41")
$s = $c.Style
$c.Value = "'-1.38%"
$c.Style = $s

$c = $ws.Range("D42")
$s = $c.Style
$c.Value = "'0.1398"
$c.Style = $s

$c = $ws.Range("E42")
$s = $c.Style
$c.Value = "'0.44%"
$c.Style = $s

$c = $ws.Range("D43")
$s = $c.Style
$c.Value = "'0.007738"
$c.Style = $s

$c = $ws.Range("E43")
$s = $c.Style
$c.Value = "'1.03%"
$c.Style = $s

$c = $ws.Range("D44")
$s = $c.Style
$c.Value = "'0.002231"
$c.Style = $s

$c = $ws.Range("E44")
$s = $c.Style
$c.Value = "'0.67%"
$c.Style = $s

$c = $ws.Range("D45")
$s = $c.Style
$c.Value = "'0.01101"
$c.Style = $s

$c = $ws.Range("E45")
$s = $c.Style
$c.Value = "'7.91%"
$c.Style = $s

$c = $ws.Range("D46")
$s = $c.Style
$c.Value = "'0.00006395"
$c.Style = $s

$c = $ws.Range("E46")
$s = $c.Style
$c.Value = "'6.99%"
$c.Style = $s

$c = $ws.Range("E47")
$s = $c.Style
$c.Value = "'0.22%"
$c.Style = $s

$c = $ws.Range("D48")
$s = $c.Style
$c.Value = "'0.0005802"
$c.Style = $s

$c = $ws.Range("E48")
$s = $c.Style
$c.Value = "'0.02%"
$c.Style = $s

$c = $ws.Range("D49")
$s = $c.Style
$c.Value = "'24.52"
$c.Style = $s

$c = $ws.Range("E49")
$s = $c.Style
$c.Value = "'182.80%"
$c.Style = $s

$c = $ws.Range("E51")
$s = $c.Style
$c.Value = "'0.22%"
$c.Style = $s
